{"js": "// 1) Remove the stray _GoBack bookmark from its original location (an\n//    otherwise-empty paragraph early in the document).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Locate the empty paragraph right after the \"2.12.2016\" entry's text\n//    (it currently has no runs) and turn it into the new \"3.12.2016\"\n//    heading paragraph, then add a following paragraph with the diary\n//    entry text (three runs) and re-create the _GoBack bookmark at its\n//    new location (end of the newly added paragraph).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet dayIdx = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"2.12.2016\") {\n    dayIdx = i;\n    break;\n  }\n}\nif (dayIdx === -1) {\n  throw new Error('Could not find the \"2.12.2016\" entry heading.');\n}\n\nlet targetIdx = -1;\nfor (let i = dayIdx + 1; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"\") {\n    targetIdx = i;\n    break;\n  }\n}\nif (targetIdx === -1) {\n  throw new Error(\"Could not find the blank paragraph to convert into the 3.12.2016 entry.\");\n}\n\nconst target = paragraphs.items[targetIdx];\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n  <w:pPr>\n    <w:jc w:val=\"center\"/>\n    <w:rPr>\n      <w:b/>\n      <w:sz w:val=\"24\"/>\n      <w:szCs w:val=\"24\"/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:b/>\n      <w:sz w:val=\"24\"/>\n      <w:szCs w:val=\"24\"/>\n    </w:rPr>\n    <w:t>3.12.2016</w:t>\n  </w:r>\n</w:p>\n<w:p>\n  <w:pPr>\n    <w:jc w:val=\"center\"/>\n    <w:rPr>\n      <w:sz w:val=\"24\"/>\n      <w:szCs w:val=\"24\"/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"24\"/>\n      <w:szCs w:val=\"24\"/>\n    </w:rPr>\n    <w:t>B\u00e6tt var vi\u00f0 bara nokkrum flash leikjum \u00ed gagnagrunninn til a\u00f0 gera hann a\u00f0eins meira dj\u00fas\u00ed, listann \u00fear a\u00f0 segja \u00feannig \u00fea\u00f0 s\u00e9 ekki bara 2 leikir \u00feegar \u00feetta er sko\u00f0a\u00f0, einig var laga\u00f0 bug a\u00f0 programmers g\u00e1tu bara sett inn einn leik \u00ed einu \u00fear sem takkinn hoppa\u00f0i ekki tilbaka \u00ed Add new game, en \u00fea\u00f0 var einf\u00f6ld lausn...</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"24\"/>\n      <w:szCs w:val=\"24\"/>\n    </w:rPr>\n    <w:t>, endurra\u00f0a\u00f0 hvernig leikirnir birtast \u00ed datagridviewinu, \u00f3\u00fearfi a\u00f0 hafa path \u00e1 undan description</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"24\"/>\n      <w:szCs w:val=\"24\"/>\n    </w:rPr>\n    <w:t>,\u00fea\u00f0 sem \u00e1 eftir er a\u00f0 gera \u00e1 einfaldann h\u00e1tt SP sem birtir \u00ed raun genre og \u00fea\u00f0 \u00e1 eftir a\u00f0 tengja genre vi\u00f0 \u00ed gagnagrunninum vi\u00f0 leikina, \u00fea\u00f0 er svosem ekki \u00ed forgangi \u00fear sem vi\u00f0 erum me\u00f0 description.</w:t>\n  </w:r>\n  <w:bookmarkStart w:id=\"1\" w:name=\"_GoBack\"/>\n  <w:bookmarkEnd w:id=\"1\"/>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the stray _GoBack bookmark from its original location (an\n#    otherwise-empty paragraph early in the document).\n$d.Bookmarks(\"_GoBack\").Delete()\n\n# 2) The empty paragraph right after the \"2.12.2016\" entry's text (#28,\n#    1-based) becomes the new \"3.12.2016\" heading paragraph, followed by\n#    a new paragraph holding the diary entry text (three runs) and the\n#    _GoBack bookmark re-created at its new location.\n$p = $d.Paragraphs(28)\n$r = $p.Range\n$r.Collapse(0)\n\n$ooxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n  <w:pPr>\n    <w:jc w:val=\"center\"/>\n    <w:rPr>\n      <w:b/>\n      <w:sz w:val=\"24\"/>\n      <w:szCs w:val=\"24\"/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:b/>\n      <w:sz w:val=\"24\"/>\n      <w:szCs w:val=\"24\"/>\n    </w:rPr>\n    <w:t>3.12.2016</w:t>\n  </w:r>\n</w:p>\n<w:p>\n  <w:pPr>\n    <w:jc w:val=\"center\"/>\n    <w:rPr>\n      <w:sz w:val=\"24\"/>\n      <w:szCs w:val=\"24\"/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"24\"/>\n      <w:szCs w:val=\"24\"/>\n    </w:rPr>\n    <w:t>B\u00e6tt var vi\u00f0 bara nokkrum flash leikjum \u00ed gagnagrunninn til a\u00f0 gera hann a\u00f0eins meira dj\u00fas\u00ed, listann \u00fear a\u00f0 segja \u00feannig \u00fea\u00f0 s\u00e9 ekki bara 2 leikir \u00feegar \u00feetta er sko\u00f0a\u00f0, einig var laga\u00f0 bug a\u00f0 programmers g\u00e1tu bara sett inn einn leik \u00ed einu \u00fear sem takkinn hoppa\u00f0i ekki tilbaka \u00ed Add new game, en \u00fea\u00f0 var einf\u00f6ld lausn...</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"24\"/>\n      <w:szCs w:val=\"24\"/>\n    </w:rPr>\n    <w:t>, endurra\u00f0a\u00f0 hvernig leikirnir birtast \u00ed datagridviewinu, \u00f3\u00fearfi a\u00f0 hafa path \u00e1 undan description</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"24\"/>\n      <w:szCs w:val=\"24\"/>\n    </w:rPr>\n    <w:t>,\u00fea\u00f0 sem \u00e1 eftir er a\u00f0 gera \u00e1 einfaldann h\u00e1tt SP sem birtir \u00ed raun genre og \u00fea\u00f0 \u00e1 eftir a\u00f0 tengja genre vi\u00f0 \u00ed gagnagrunninum vi\u00f0 leikina, \u00fea\u00f0 er svosem ekki \u00ed forgangi \u00fear sem vi\u00f0 erum me\u00f0 description.</w:t>\n  </w:r>\n  <w:bookmarkStart w:id=\"1\" w:name=\"_GoBack\"/>\n  <w:bookmarkEnd w:id=\"1\"/>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n\n$r.InsertXML($ooxml)\n"}
